$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# A new "Main Menu" task is being added to the Menus/UI/UX breakdown list
# (column C, rows 12-17). This requires inserting a new cell at C13 and
# shifting the existing C13:C17 values down into C14:C18.
#
# Likewise, a blank separator cell is inserted into column C at row 22 (to
# line up with a new blank row above the "Asset collection" section), which
# shifts the existing C22:C86 values down into C23:C87.
#
# Because this runtime's Range.Insert(xlShiftDown) shifts the *whole row*
# (all columns) instead of Excel's real behaviour of shifting only the
# selected column, we emulate a column-only "insert/shift down" by manually
# copying each source cell's content/formula/font down into the cell below,
# working from the bottom of the range upwards so nothing is overwritten
# before it is read.
# ---------------------------------------------------------------------------

function Shift-ColumnCDown($TopRow, $BottomRow) {
    for ($r = $BottomRow; $r -ge $TopRow; $r--) {
        $src = $ws.Cells.Item($r, 3)
        $dst = $ws.Cells.Item($r + 1, 3)

        $srcHasFormula = $src.HasFormula
        $srcVal = $src.Value2
        $srcBold = $src.Font.Bold
        $srcUnderline = $src.Font.Underline

        # Wipe the destination completely first (value + formatting) so
        # cells that should end up blank don't retain stray formatting.
        $dst.Clear()

        if ($srcHasFormula) {
            $dst.Formula = $src.Formula
        } elseif ($srcVal -ne $null -and $srcVal -ne "") {
            $dst.Value = $srcVal
        }

        if ($srcBold -eq $true) {
            $dst.Font.Bold = $true
        }
        if ($srcUnderline -eq 2) {
            $dst.Font.Underline = 2
        }
    }
}

# --- Block 1: rows 12-17 -> 12-18, new row 13 = "Main Menu" ---------------
Shift-ColumnCDown 13 17
$ws.Range("C13").Clear()
$ws.Range("C13").Value = "Main Menu"

# --- Block 2: rows 22-86 -> 23-87, new row 22 left blank -------------------
Shift-ColumnCDown 22 86
$ws.Range("C22").Clear()

# Match the active selection left behind by the author's edit.
$ws.Range("C12").Select()
